$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.381.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.882.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.36%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.7130"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.01%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'242.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.28%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value = "'0.08032"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.65%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.85%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'25.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.26%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08351"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.86%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.885.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.24%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.75%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.7192"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.27%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'94.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.20%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.334"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +5.47%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008551"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +4.25%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'29.385.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.25%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'BitcoinCash"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'242.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.06%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'2.135.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.02%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.14%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.12%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.74%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.11%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.84%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'163.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.38%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.086"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.38%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'18.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.76%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.511"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.19%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.422"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.58%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.333"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.38%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -6.55%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05390"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.52%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.951"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.85%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.73%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7512"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.14%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.695"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.32%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.41%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.287.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +8.83%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.06%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.602"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.43%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.9175"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.34%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'74.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +2.48%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'111.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +5.15%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.09%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +6.30%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.047.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.85%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.810"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.05%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.5221"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.28%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'9.549"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.86%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.4393"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.89%  "
$ws.Range("E51").Style = "Normal"
